$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "info"
$ws2 = $wb.Worksheets.Item(2)   # "GWP"

# ---------------------------------------------------------------------------
# 1) Sheet "info": add a "kind" column, and two new rows (Electricity, NaClO)
#    Order of writes matters so that new shared-string entries are appended
#    in the same sequence as the target workbook.
# ---------------------------------------------------------------------------
$ws1.Range("A12").Value = "Electricity"
$ws1.Range("B12").Value = "kWh"
$ws1.Range("C12").Value = "ImpactItem"

$ws1.Range("C13").Value = "StreamImpactItem"
$ws1.Range("A13").Value = "NaClO"
$ws1.Range("B13").Value = "kg"

$ws1.Range("C2:C11").Value = "ImpactItem"

$ws1.Range("C1").Value = "kind"
$ws1.Range("C1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) Sheet "GWP": append impact-factor rows for Electricity and NaClO
# ---------------------------------------------------------------------------
$ws2.Range("A12").Value = "Electricity"
$ws2.Range("B12").Value = "kg CO2-eq"
$ws2.Range("C12").Value = 0.1135
$ws2.Range("D12").Value = 0.106
$ws2.Range("E12").Value = 0.121
$ws2.Range("F12").Value = "uniform"
$ws2.Range("G12").Value = "ecoinvent 3"

$ws2.Range("A13").Value = "NaClO"
$ws2.Range("B13").Value = "kg CO2-eq"
$ws2.Range("C13").Value = 2.6287
$ws2.Range("D13").Formula = '=$C13*0.75'
$ws2.Range("E13").Formula = '=$C13*1.25'
$ws2.Range("F13").Value = "uniform"
$ws2.Range("G13").Value = "ecoinvent 3"

# ---------------------------------------------------------------------------
# 3) Styling: drop the blue highlight fill that used to mark the numeric
#    columns and keep the header row bold. The two-step
#    (set a color, then clear it) sequence below forces Excel to persist an
#    explicit "no fill" flag on every touched cell, matching how the
#    original author removed the fill via the ribbon.
# ---------------------------------------------------------------------------
$ws2.Range("A1:G1").Interior.ColorIndex = 5
$ws2.Range("A1:G1").Interior.Pattern = -4142
$ws2.Range("A1:G1").Font.Bold = $true

$ws2.Range("A2:G12").Interior.ColorIndex = 5
$ws2.Range("A2:G12").Interior.Pattern = -4142

$ws2.Range("B13").Interior.ColorIndex = 5
$ws2.Range("B13").Interior.Pattern = -4142
$ws2.Range("F13:G13").Interior.ColorIndex = 5
$ws2.Range("F13:G13").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 4) View state: make "GWP" the active/selected sheet, update selections
# ---------------------------------------------------------------------------
$ws1.Range("A14:XFD14").Select()
$ws2.Range("F17").Select()
$ws2.Activate()
